# Update "Taetigkeitsprotokoll" worksheet with new activity entries
# (commit: "Update Kurs ohne Where!!!!")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 10 and 11 descriptions ---
# Row 10 (26.11.2019): "Programming" -> longer description (wraps to 2 lines)
$ws.Range("E10").Value = "stylesheet, Registrierung, Password-Class, Logout, Weiterleitung Registrierung"

# Row 11 (03.12.2019): "Programming" -> longer description
$ws.Range("E11").Value = "Masterunterscheidung User/Admin, Passwort-Überprüfung, E-Mail-Klasse, AGBs und DsGV, Logout-Button, Weiterleitung Passwort vergessen"

# Row 10 now wraps onto two lines, bump its row height to match
$ws.Rows.Item(10).RowHeight = 28.5

# --- Fill in new activity rows 12-15 ---
# Row 12: Di 10.12.2019
$ws.Range("A12").Value = "Di"
$ws.Range("B12").Value = "10.12.2019"
$ws.Range("C12").Value = 0.32291666666666669
$ws.Range("D12").Value = 0.54861111111111105
$ws.Range("E12").Value = "Kurse anzeigen+anlegen"

# Row 13: Mi 11.12.2019
$ws.Range("A13").Value = "Mi"
$ws.Range("B13").Value = "11.12.2019"
$ws.Range("C13").Value = 0.40277777777777773
$ws.Range("D13").Value = 0.47222222222222227
$ws.Range("E13").Value = "Kurs anlegen"

# Row 14: Di 17.12.2019
$ws.Range("A14").Value = "Di"
$ws.Range("B14").Value = "17.12.2019"
$ws.Range("C14").Value = 0.32291666666666669
$ws.Range("D14").Value = 0.54861111111111105
$ws.Range("E14").Value = "Kurs anlegen"

# Row 15: Di 7.1.2020
$ws.Range("A15").Value = "Di"
$ws.Range("B15").Value = "7.1.2020"
$ws.Range("C15").Value = 0.32291666666666669
$ws.Range("D15").Value = 0.54861111111111105
$ws.Range("E15").Value = "Kurs anlegen"

# --- Update the selected cell to match final cursor position ---
$ws.Range("E19").Select() | Out-Null
